$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (old C/D/E shift right to D/E/F)
[void]$ws.Columns("C").Insert()

# Match the new column's width to column B's width (24)
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Header for the new column
$ws.Range("C1").Value = "BaseDialog2"

# New base-dialog values for the vehicle-insurance rows
$ws.Range("C8").Value = "dlgProductInsuranceBase"
$ws.Range("C9").Value = "dlgProductInsuranceBase"
$ws.Range("C10").Value = "dlgProductInsuranceBase"
$ws.Range("C11").Value = "dlgProductInsuranceBase"

# Update the active selection
[void]$ws.Range("C12").Select()

# Try to mirror the stored window geometry (best effort)
$win = $wb.Windows.Item(1)
$win.Left = 4572
$win.Top = 3024
$win.Width = 22560
$win.Height = 12480
